# Fix typo on slide 3: remove the duplicated "concatenation, " that was
# accidentally repeated in the list of allowed Matlab operations.
#
# Old: "concatenation, matrix multiplication, concatenation, transpose and inv "
# New: "concatenation, matrix multiplication, transpose and inv "
#
# The affected run keeps its existing formatting (blue font color), so we
# replace only the offending substring in-place via TextRange.Characters()
# instead of overwriting the whole TextRange (which would lose per-run
# formatting boundaries).

$old = "concatenation, matrix multiplication, concatenation, transpose and inv "
$new = "concatenation, matrix multiplication, transpose and inv "

$p = $ppt.ActivePresentation

for ($slideIdx = 1; $slideIdx -le $p.Slides.Count; $slideIdx++) {
    $slide = $p.Slides.Item($slideIdx)
    $shapes = $slide.Shapes

    for ($shapeIdx = 1; $shapeIdx -le $shapes.Count; $shapeIdx++) {
        $shape = $shapes.Item($shapeIdx)

        if ($shape.HasTextFrame) {
            $textRange = $shape.TextFrame.TextRange
            $fullText = $textRange.Text

            $pos = $fullText.IndexOf($old)
            while ($pos -ge 0) {
                $targetRange = $textRange.Characters($pos + 1, $old.Length)
                $targetRange.Text = $new

                $fullText = $textRange.Text
                $pos = $fullText.IndexOf($old)
            }
        }
    }
}
